# "Day 1, designed layout"
#
# The Gantt-chart template had a handful of label/blank cells (A3, A5, A7,
# A13 and the B15:D15 footer labels) carrying a cell format that was an
# exact duplicate of the plain bordered style already used by the
# surrounding blank cells (e.g. A4). Re-apply that shared format so those
# cells stop referencing the redundant, duplicate style record - this is
# what "designing the layout" collapsed down to once the sheet was
# reopened/resaved.
#
# A9 and A11 (the wrap-text task rows) keep their wrapped, bordered
# look untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the already-common plain bordered format (as seen on A4) as the
# canonical style and copy it onto the cells that were carrying the
# redundant duplicate.
$ws.Range("A4").Copy()

$labelCells = @("A3", "A5", "A7", "A13", "B15", "C15", "D15")
foreach ($cellRef in $labelCells) {
    $ws.Range($cellRef).PasteSpecial(-4122)  # xlPasteFormats
}

$excel.CutCopyMode = $false

# Leave the cursor where the file was last saved.
$ws.Range("G11").Select()
